$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'59.459.61"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = "'2.636.99"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'536.03"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.17%  '
$ws.Range('D6').Value = "'144.57"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.73%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('D9').Value = "'2.648.58"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('D10').Value = "'6.63"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E11').Value = '  -1.10%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('D14').Value = "'3.108.98"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.27%  '
$ws.Range('D15').Value = "'59.362.44"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').Value = "'21.06"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.21%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = "'2.661.07"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.89%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = "'0.0000134"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = "'339.25"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('D20').Value = "'4.38"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').Value = "'10.32"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.74%  '
$ws.Range('D22').Value = "'6.29"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('D23').Value = "'0.999"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = "'67.03"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = "'0.414"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('D26').Value = "'0.164"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.40%  '
$ws.Range('D27').Value = "'0.999"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').Value = "'7.25"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('D29').Value = "'0.0₃0742"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.37%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').Value = "'1.64"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('D33').Value = "'18.86"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').Value = "'151.40"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.71%  '
$ws.Range('D35').Value = "'3.98"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('E36').Value = '  +0.79%  '
$ws.Range('D37').Value = "'0.841"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('D38').Value = "'0.831"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.66%  '
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').Value = "'288.12"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.19%  '
$ws.Range('D41').Value = "'3.58"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = "'0.602"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.61%  '
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('D45').Value = "'19.27"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.09%  '
$ws.Range('D46').Value = "'0.0532"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('D47').Value = "'0.0945"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.78%  '
$ws.Range('D48').Value = "'1.968.32"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('D50').Value = "'4.55"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('D51').Value = "'18.23"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.64%  '
